$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the volatile RAND() driver cell (keep its style)
$ws.Range("V4").ClearContents()

# Remove the helper index row (U5:AJ5) entirely
$ws.Range("U5:AJ5").ClearContents()

# Clear the per-subject control-threshold values (U8:AJ8), keep formatting
$ws.Range("U8:AJ8").ClearContents()

# Clear the dependent formulas (U9:AJ9), keep formatting
$ws.Range("U9:AJ9").ClearContents()

# Update the view: drop the scrolled-to-T1 position and move the selection
$ws.Range("J32").Select()
